$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "Dr. Eman Tantawi, Dr. Servinaz Sayed Mohammad, Dr. Hend Mahmoud, Dr. Majorelle Magdy"
$ws.Range("G3").Value = "Dr. Menna tuâ€™Allah Medhat, Dr. Amira Sobhy, Dr. Eman Tantawi, Dr. Veronia Rafat, Dr. Asmaa Reda"
$ws.Range("G4").Value = "Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Eman Tantawi, Dr. Veronia Rafat, Dr. Rana Abo-Zaid, Dr. Asmaa Reda"
$ws.Range("G5").Value = "Dr. Hend Mahmoud, Dr. Eman Tantawi, Dr. Veronia Rafat, Dr. Nesma, Dr. Servinaz Sayed Mohammad, Dr. Mohammad El-Tanany, Dr. Nourhan Mahmoud, Dr. Hanan Ragab"
$ws.Range("G6").Value = "Dr. Menna tuâ€™Allah Medhat, Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Eman Tantawi, Dr. Veronia Rafat, Dr. Servinaz Sayed Mohammad, Dr. Nahla Nagiub, Dr. Nourhan Mahmoud, Dr. Gehan Adel, Dr. Asmaa Reda"
$ws.Range("G7").Value = "Dr. Menna tuâ€™Allah Medhat, Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Eman Tantawi, Dr. Veronia Rafat, Dr. Rana Abo-Zaid, Dr. Servinaz Sayed Mohammad, Dr. Gehan Adel, Dr. Asmaa Reda"
$ws.Range("G8").Value = "Dr. Majorelle Magdy, Dr. Eman Tantawi, Administrator, Dr. Shimaa Ahmad Mekki, Dr. Manar Montaser, Dr. Asmaa Reda"
$ws.Range("G9").Value = "Dr. Menna tuâ€™Allah Medhat, Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Majorelle Magdy, Dr. Gehan Adel, Dr. Rana Abo-Zaid, Dr. Manar Montaser, Dr. Asmaa Reda"
$ws.Range("G10").Value = "Dr. Alshimaa Atef, Dr. Sara Wael, Dr. Rana Abo-Zaid, Dr. Servinaz Sayed Mohammad, Dr. Shimaa Ahmad Mekki, Dr. Heba Mahmoud Ali, Dr. Gehan Adel"
$ws.Range("G11").Value = "Dr. Eman Tantawi, Dr. Veronia Rafat, Dr. Asmaa Reda, Dr. Hend Mahmoud"
$ws.Range("G13").Value = "Dr. Mariam Nour El-Din, D Wessam Atef, Dr. Safa Hany, Dr. Omnia Mohammad, Dr. Shimaa Ashraf"
$ws.Range("G16").Value = "Dr. Amal Awwad, Dr. Nourhan Mohammad"
$ws.Range("G17").Value = "Dr. Nourhan Osama, Dr. Marwa Mustafa, Dr. Yasmeena Fattoh, Dr. Sarah Abdelmohsen, Dr. Madeha Saeed, Dr. Arwa Al-Sayed, Dr. Basma Hamed, Dr. Dina Adel, Dr. Eman M. Abo-Sakaya, Dr. Esraa Mostafa"
$ws.Range("G19").Value = "D Mariam E. Mohammad, Dr. Sarah Mahdy"
$ws.Range("G23").Value = "Dr. Hana Amr, Dr. Nourham Mostafa"
$ws.Range("G24").Value = "Dr. Youstina Magdy, Dr. Wafaa Ebida, Dr. Aya Emad, Dr. Salma Hassan, Dr. Marina Atef, Dr. Neveen Nashaat, Dr. Remon, Dr. Maryam Ashraf, Dr. Yasmin, Dr. Monica, Dr. Ola Abd Al-Fattah"
$ws.Range("G25").Value = "Dr. Abdullah El-Agrody, Dr. Aya Emad, Dr. Marina Atef, Dr. Remon, Dr. Youstina Magdy, Dr. Eman Samir Gabry, Dr. Ola Abd Al-Fattah"
$ws.Range("G26").Value = "Dr. Youstina Magdy, Dr. Gehad Salah"
$ws.Range("G27").Value = "Dr. Wafaa Ebida, Dr. Salma Hassan, Dr. Remon, Dr. Eman Mohammad Al, Dr. Yasmin, Dr. Neveen Nashaat, Dr. Eman Samir Gabry, Dr. Ola Abd Al-Fattah"
$ws.Range("G28").Value = "Dr. Abdullah El-Agrody, Dr. Nardine, Dr. Wafaa Ebida, Dr. Aya Hanafy, Dr. Salma Hassan, Dr. Remon, Dr. Neveen Nashaat, Dr. Eman Samir Gabry"
$ws.Range("G29").Value = "Dr. Remon, Dr. Naema Gomaa, Dr. Neveen Nashaat, Dr. Monica, Dr. Eman Samir Gabry, Dr. Ola Abd Al-Fattah"
$ws.Range("G30").Value = "Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Eman Tantawi, Dr. Veronia Rafat, Dr. Rana Abo-Zaid, Dr. Asmaa Reda"
$ws.Range("G31").Value = "Dr. Menna tuâ€™Allah Medhat, Dr. Amira Sobhy, Dr. Eman Tantawi, Dr. Veronia Rafat, Dr. Asmaa Reda"
$ws.Range("G32").Value = "Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Eman Tantawi, Dr. Veronia Rafat, Dr. Rana Abo-Zaid, Dr. Asmaa Reda"
$ws.Range("G33").Value = "Dr. Hend Mahmoud, Dr. Eman Tantawi, Dr. Veronia Rafat, Dr. Nesma, Dr. Servinaz Sayed Mohammad, Dr. Mohammad El-Tanany, Dr. Nourhan Mahmoud, Dr. Hanan Ragab"
$ws.Range("G34").Value = "Dr. Menna tuâ€™Allah Medhat, Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Eman Tantawi, Dr. Veronia Rafat, Dr. Servinaz Sayed Mohammad, Dr. Nahla Nagiub, Dr. Nourhan Mahmoud, Dr. Gehan Adel, Dr. Asmaa Reda"
$ws.Range("G35").Value = "Dr. Menna tuâ€™Allah Medhat, Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Eman Tantawi, Dr. Veronia Rafat, Dr. Rana Abo-Zaid, Dr. Servinaz Sayed Mohammad, Dr. Gehan Adel, Dr. Asmaa Reda"
$ws.Range("G36").Value = "Dr. Majorelle Magdy, Dr. Eman Tantawi, Administrator, Dr. Shimaa Ahmad Mekki, Dr. Manar Montaser, Dr. Asmaa Reda"
$ws.Range("G37").Value = "Dr. Menna tuâ€™Allah Medhat, Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Majorelle Magdy, Dr. Gehan Adel, Dr. Rana Abo-Zaid, Dr. Manar Montaser, Dr. Asmaa Reda"
$ws.Range("G38").Value = "Dr. Alshimaa Atef, Dr. Sara Wael, Dr. Rana Abo-Zaid, Dr. Servinaz Sayed Mohammad, Dr. Shimaa Ahmad Mekki, Dr. Heba Mahmoud Ali, Dr. Gehan Adel"
$ws.Range("G39").Value = "Dr. Eman Tantawi, Dr. Veronia Rafat, Dr. Asmaa Reda, Dr. Hend Mahmoud"
$ws.Range("G41").Value = "Dr. Mariam Nour El-Din, D Wessam Atef, Dr. Safa Hany, Dr. Omnia Mohammad, Dr. Shimaa Ashraf"
$ws.Range("G44").Value = "Dr. Amal Awwad, Dr. Nourhan Mohammad"
$ws.Range("G45").Value = "Dr. Nourhan Osama, Dr. Marwa Mustafa, Dr. Yasmeena Fattoh, Dr. Sarah Abdelmohsen, Dr. Madeha Saeed, Dr. Arwa Al-Sayed, Dr. Basma Hamed, Dr. Dina Adel, Dr. Eman M. Abo-Sakaya, Dr. Esraa Mostafa"
$ws.Range("G47").Value = "D Mariam E. Mohammad, Dr. Sarah Mahdy"
$ws.Range("G51").Value = "Dr. Hana Amr, Dr. Nourham Mostafa"
$ws.Range("G52").Value = "Dr. Youstina Magdy, Dr. Wafaa Ebida, Dr. Aya Emad, Dr. Salma Hassan, Dr. Marina Atef, Dr. Neveen Nashaat, Dr. Remon, Dr. Maryam Ashraf, Dr. Yasmin, Dr. Monica, Dr. Ola Abd Al-Fattah"
$ws.Range("G53").Value = "Dr. Abdullah El-Agrody, Dr. Aya Emad, Dr. Marina Atef, Dr. Remon, Dr. Youstina Magdy, Dr. Eman Samir Gabry, Dr. Ola Abd Al-Fattah"
$ws.Range("G54").Value = "Dr. Youstina Magdy, Dr. Gehad Salah"
$ws.Range("G55").Value = "Dr. Wafaa Ebida, Dr. Salma Hassan, Dr. Remon, Dr. Eman Mohammad Al, Dr. Yasmin, Dr. Neveen Nashaat, Dr. Eman Samir Gabry, Dr. Ola Abd Al-Fattah"
$ws.Range("G56").Value = "Dr. Abdullah El-Agrody, Dr. Nardine, Dr. Wafaa Ebida, Dr. Aya Hanafy, Dr. Salma Hassan, Dr. Remon, Dr. Neveen Nashaat, Dr. Eman Samir Gabry"
$ws.Range("G57").Value = "Dr. Remon, Dr. Naema Gomaa, Dr. Neveen Nashaat, Dr. Monica, Dr. Eman Samir Gabry, Dr. Ola Abd Al-Fattah"
